$d = $word.ActiveDocument

# Locate the paragraph that credits the map author / links to the astromap,
# regardless of its numeric index in the document.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Jenik Hollan*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    Write-Host "Target paragraph not found"
} else {
    $full = $target.Range
    # Exclude the trailing paragraph mark from the range so InsertXML only
    # rewrites the runs, leaving the paragraph's own pPr/formatting intact.
    $r = $d.Range($full.Start, $full.End - 1)

    $newText = "Peta di dokumen ini disiapkan oleh Jenik Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p><w:r/><w:r><w:t>' + $newText + '</w:t></w:r></w:p></w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $r.InsertXML($xml)
    Write-Host "Updated astromap year in credit paragraph."
}
